$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing review (row 2) down to row 3, then add two new rows ---

# Remember the original row 2 values before we overwrite anything (".Value"
# getter is unreliable in this host - use ".Value2" for reads).
$origA2 = $ws.Range("A2").Value2
$origB2 = $ws.Range("B2").Value2
$origC2 = $ws.Range("C2").Value2
$origD2 = $ws.Range("D2").Value2
$origE2 = $ws.Range("E2").Value2
$origF2 = $ws.Range("F2").Value2

# Stash a copy of row 2's plain (non-hyperlinked) cell formats off to the
# side before anything touches styles, so we can restore the C/D look after
# Hyperlinks.Add later stamps its own "Hyperlink" style onto those cells.
$ws.Range("C2:D2").Copy()
$ws.Range("H2:I2").PasteSpecial(-4122)

# Copy row 2's formatting down into rows 3 and 4 before changing any values.
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)
$ws.Range("A2:F2").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Drop the old hyperlinks (this clears every hyperlink on the sheet in this
# engine, so we recreate the full set below in the final desired layout).
$ws.Range("C2").Hyperlinks.Delete()

# Row 3 = the original review, moved down one row.
$ws.Range("A3").Value = $origA2
$ws.Range("B3").Value = $origB2
$ws.Range("C3").Value = $origC2
$ws.Range("D3").Value = $origD2
$ws.Range("E3").Value = $origE2
$ws.Range("F3").Value = $origF2

# Row 2 = new "taxi game" review.
$ws.Range("A2").Value = "com.singleton.stretchy"
$ws.Range("B2").Value = "taxi game"
$ws.Range("C2").Value = "budoyoni@gmail.com"
$ws.Range("D2").Value = "sm6502345@gmail.com"
$ws.Range("E2").Value = "27/5/2019 15:59"
$ws.Range("F2").Value = "wow, this cars game is so good. I think it is fantastic and fun. The levels are so unique and graphics. Love it!"

# Row 4 = new second bitcoin-guide review.
$ws.Range("A4").Value = "com.hamxa.shaynachim"
$ws.Range("B4").Value = "bitcoin guide"
$ws.Range("C4").Value = "zaittomer@gmail.com "
$ws.Range("D4").Value = "eligitel@gmail.com"
$ws.Range("E4").Value = "27/5/2019 15:59"
$ws.Range("F4").Value = "Exclusive info and great explanations!! bitcoin is hottt"

# Recreate hyperlinks for every email cell, in final row order, so the
# relationship ids line up the way they would after a real insert.
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:budoyoni@gmail.com", "", "", "budoyoni@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:sm6502345@gmail.com", "", "", "sm6502345@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:eligitel@gmail.com", "", "", "eligitel@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:ronenchen27@gmail.com", "", "", "ronenchen27@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:zaittomer@gmail.com ", "", "", "zaittomer@gmail.com ")

# Hyperlinks.Add stamps the built-in "Hyperlink" style onto the cell; restore
# the plain formatting (style index 2 in the original file) from the stash.
$ws.Range("H2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("I2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clear the scratch cells used to stash formatting.
$ws.Range("H2:I2").Clear()

# Row 4 gets a slightly taller row height, matching the source edit.
$ws.Rows.Item(4).RowHeight = 13.8

$ws.Range("F4").Select()
